# Adds five new question sheets (11_ .. 15_) to the workbook, matching the
# "Add files via upload" commit: three sheets of new multiple-choice
# questions (low-pass filters, circuits/electrons, voltage) plus two blank
# placeholder sheets.

$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd($wb, $name) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $last)
    $ws.Name = $name
    return $ws
}

# Helper to style a "question" cell (column A data rows): left/top, wrapped.
function Set-QuestionStyle($range) {
    $range.HorizontalAlignment = -4131   # xlLeft
    $range.VerticalAlignment = -4160     # xlTop
    $range.WrapText = $true
}

# ---------------------------------------------------------------------
# Sheet "11_" - low pass filter question
# ---------------------------------------------------------------------
$ws11 = Add-SheetAtEnd $wb "11_"

$ws11.Columns.Item(1).ColumnWidth = 27.28515625
$ws11.Columns.Item(2).ColumnWidth = 23.85546875
$ws11.Columns.Item(3).ColumnWidth = 51.85546875

# Column A default alignment (left/top, no wrap) matches the workbook's
# column-level style for this sheet.
$colA = $ws11.Columns.Item(1)
$colA.HorizontalAlignment = -4131
$colA.VerticalAlignment = -4160

$data11 = @(
    @("What does a low pass filter do?", "Correct", "Comment"),
    @("It reduces the amplitude of an incoming signal", "N", ""),
    @("It removes the low frequency part of an incoming signal", "N", ""),
    @("It removes the high frequency part of an incoming signal", "Y", "Yep!   This is why it is called ""low pass"": low frequencies are allowed to ""pass"", while high frequencies are removed."),
    @("It removes out the high amplitude part of an incoming signal", "N", ""),
    @("It removes ""noise"" from an incoming signal", "N", "This is often the purpose of a low pass filter, but ""noise"" in a signal can be of any frequency: ""noise"" is simply defined as the part of the signal that we are not interested in, just like a ""weed"" is just a plant we don't want in the garden.   So if the ""noise"" is a high frequency oscillations, then a low pass filter would remove it; but if the ""noise"" is a low frequency hum, it would not.")
)

for ($i = 0; $i -lt $data11.Length; $i++) {
    $row = $i + 1
    $ws11.Cells.Item($row, 1).Value = $data11[$i][0]
    $ws11.Cells.Item($row, 2).Value = $data11[$i][1]
    if ($data11[$i][2] -ne "") {
        $ws11.Cells.Item($row, 3).Value = $data11[$i][2]
    }
    Set-QuestionStyle $ws11.Cells.Item($row, 1)
}

$ws11.Range("A1:C6").WrapText = $true

# ---------------------------------------------------------------------
# Sheet "12_" - what drives electrons question
# ---------------------------------------------------------------------
$ws12 = Add-SheetAtEnd $wb "12_"

$ws12.Columns.Item(1).ColumnWidth = 25.28515625
$ws12.Columns.Item(3).ColumnWidth = 36.140625

$data12 = @(
    @("What drives electrons to move around a circuit?", ""),
    @("Current ", "N"),
    @("Voltage difference", "Y", "Yep!  Electrons want to move from a high voltage area to a low voltage area. "),
    @("Resistance", "N"),
    @("Capacitance", "N")
)

for ($i = 0; $i -lt $data12.Length; $i++) {
    $row = $i + 1
    $ws12.Cells.Item($row, 1).Value = $data12[$i][0]
    if ($data12[$i].Length -gt 1 -and $data12[$i][1] -ne "") {
        $ws12.Cells.Item($row, 2).Value = $data12[$i][1]
    }
    if ($data12[$i].Length -gt 2 -and $data12[$i][2] -ne "") {
        $ws12.Cells.Item($row, 3).Value = $data12[$i][2]
    }
}

$ws12.Range("A1:G10").WrapText = $true

# ---------------------------------------------------------------------
# Sheet "13_" - voltage question (multi-select)
# ---------------------------------------------------------------------
$ws13 = Add-SheetAtEnd $wb "13_"

$ws13.Columns.Item(1).ColumnWidth = 34.28515625
$ws13.Columns.Item(3).ColumnWidth = 39.7109375

$data13 = @(
    @("Which of the following is a good way to think of voltage?  Mark all that are true", ""),
    @("Voltage is a kind of energy", "Y", "Yep!  Electrons don't like to be next to each other, so pushing a bunch of them together requires work.   This work gets stored as ""potential energy"" that can cause the electrons to move if they are allowed to."),
    @("Voltage is like pressure in a fluid", "Y", "Yep!  Just like pressure difference causes the force that drives fluid through a pipe, voltage difference causes the force that drives a current through a circuit"),
    @("Voltage is a place in the circuit where electrons are packed together", "Y", "Yep!  The positive terminal of a battery is loaded with electrons packed together, just rarin' to get away.   An electron that reaches the negative terminal is moved back to the positive terminal by the stored energy of the battery: that's what a battery does.   (Note: this is all by convention.  In actuality, the process happens in reverse)."),
    @("Voltage supplies the super power of any super hero with a lightning bolt on his or her chest", "Y", "Yep!  This is hard to deny.")
)

for ($i = 0; $i -lt $data13.Length; $i++) {
    $row = $i + 1
    $ws13.Cells.Item($row, 1).Value = $data13[$i][0]
    if ($data13[$i].Length -gt 1 -and $data13[$i][1] -ne "") {
        $ws13.Cells.Item($row, 2).Value = $data13[$i][1]
    }
    if ($data13[$i].Length -gt 2 -and $data13[$i][2] -ne "") {
        $ws13.Cells.Item($row, 3).Value = $data13[$i][2]
    }
}

$ws13.Range("A1:H17").WrapText = $true

# ---------------------------------------------------------------------
# Sheets "14_" and "15_" - blank placeholder sheets
# ---------------------------------------------------------------------
$ws14 = Add-SheetAtEnd $wb "14_"
$ws15 = Add-SheetAtEnd $wb "15_"
